# Weekly refresh of "Fruto del paraíso" (Vega Modelo de Temuco) price rows.
# The underlying data rows (2-18, except 8 and 14 which are unchanged) get
# re-sorted / refreshed with a new weekly snapshot: dates, quality, volume,
# prices and the commercialization-unit text are updated per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45037
$ws.Range("N2").Value = "$/caja 15 kilos empedrada"
$ws.Range("P2").Value = 1600
$ws.Range("Q2").Value = 15

# Row 3
$ws.Range("D3").Value = 44285
$ws.Range("J3").Value = 20
$ws.Range("K3").Value = 25000
$ws.Range("L3").Value = 25000
$ws.Range("M3").Value = 25000
$ws.Range("N3").Value = "$/caja 18 kilos empedrada"
$ws.Range("P3").Value = 1389
$ws.Range("Q3").Value = 18

# Row 4
$ws.Range("D4").Value = 45044
$ws.Range("J4").Value = 40
$ws.Range("N4").Value = "$/caja 18 kilos empedrada"

# Row 5
$ws.Range("D5").Value = 44293
$ws.Range("J5").Value = 10
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = 25000
$ws.Range("N5").Value = "$/caja 15 kilos empedrada"
$ws.Range("P5").Value = 1667
$ws.Range("Q5").Value = 15

# Row 6
$ws.Range("D6").Value = 44315
$ws.Range("I6").Value = "Especial"
$ws.Range("J6").Value = 10
$ws.Range("K6").Value = 30000
$ws.Range("L6").Value = 30000
$ws.Range("M6").Value = 30000
$ws.Range("N6").Value = "$/caja 20 kilos empedrada"
$ws.Range("P6").Value = 1500
$ws.Range("Q6").Value = 20

# Row 7
$ws.Range("D7").Value = 44315
$ws.Range("J7").Value = 20
$ws.Range("K7").Value = 15000
$ws.Range("L7").Value = 15000
$ws.Range("M7").Value = 15000
$ws.Range("N7").Value = "$/caja 15 kilos granel"
$ws.Range("P7").Value = 1000
$ws.Range("Q7").Value = 15

# Row 9
$ws.Range("D9").Value = 45033
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 24000
$ws.Range("L9").Value = 24000
$ws.Range("M9").Value = 24000
$ws.Range("N9").Value = "$/caja 18 kilos granel"
$ws.Range("P9").Value = 1333
$ws.Range("Q9").Value = 18

# Row 10
$ws.Range("D10").Value = 44280
$ws.Range("J10").Value = 30
$ws.Range("K10").Value = 25000
$ws.Range("L10").Value = 25000
$ws.Range("M10").Value = 25000
$ws.Range("N10").Value = "$/caja 18 kilos empedrada"
$ws.Range("P10").Value = 1389
$ws.Range("Q10").Value = 18

# Row 11
$ws.Range("D11").Value = 45034
$ws.Range("N11").Value = "$/caja 18 kilos granel"

# Row 12
$ws.Range("D12").Value = 45041
$ws.Range("J12").Value = 80

# Row 13
$ws.Range("D13").Value = 45040
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 18000
$ws.Range("L13").Value = 18000
$ws.Range("M13").Value = 18000
$ws.Range("N13").Value = "$/caja 18 kilos empedrada"
$ws.Range("P13").Value = 1000
$ws.Range("Q13").Value = 18

# Row 15
$ws.Range("D15").Value = 45015
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 50
$ws.Range("K15").Value = 24000
$ws.Range("L15").Value = 24000
$ws.Range("M15").Value = 24000
$ws.Range("N15").Value = "$/caja 18 kilos empedrada"
$ws.Range("P15").Value = 1333
$ws.Range("Q15").Value = 18

# Row 16
$ws.Range("D16").Value = 45043
$ws.Range("J16").Value = 80
$ws.Range("K16").Value = 24000
$ws.Range("L16").Value = 24000
$ws.Range("M16").Value = 24000
$ws.Range("N16").Value = "$/caja 18 kilos empedrada"
$ws.Range("P16").Value = 1333
$ws.Range("Q16").Value = 18

# Row 17
$ws.Range("D17").Value = 44313
$ws.Range("J17").Value = 20
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("M17").Value = 15000
$ws.Range("N17").Value = "$/caja 15 kilos empedrada"
$ws.Range("P17").Value = 1000
$ws.Range("Q17").Value = 15

# Row 18
$ws.Range("D18").Value = 44313
$ws.Range("J18").Value = 20
$ws.Range("K18").Value = 30000
$ws.Range("L18").Value = 30000
$ws.Range("M18").Value = 30000
$ws.Range("N18").Value = "$/caja 20 kilos empedrada"
$ws.Range("P18").Value = 1500
$ws.Range("Q18").Value = 20
